$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New device rows for regcntr_id 10002, continuing the device_id sequence
$newRows = @(
    @(10002, 3000176),
    @(10002, 3000177),
    @(10002, 3000178),
    @(10002, 3000179),
    @(10002, 3000180)
)

$startRow = 157
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $regcntrId = $newRows[$i][0]
    $deviceId = $newRows[$i][1]

    $ws.Cells.Item($r, 1).Value = $regcntrId
    $ws.Cells.Item($r, 2).Value = $deviceId
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
}

# Scroll / select to match the new view state after appending rows
$ws.Range("C158").Select()
